$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the pin names in column D (rows 4-15) per the new data sheet
$ws.Range("D4").Value = "SRAM WE"
$ws.Range("D5").Value = "SRAM OE"
$ws.Range("D6").Value = "SRAM CE"
$ws.Range("D8").Value = "ADDR <14>"
$ws.Range("D7").Value = "ADDR <15>"
$ws.Range("D9").Value = "ADDR <13>"
$ws.Range("D10").Value = "GND"
$ws.Range("D11").Value = "ADDR <12>"
$ws.Range("D12").Value = "ADDR <11>"
$ws.Range("D13").Value = "ADDR <10>"
$ws.Range("D14").Value = "ADDR <9>"
$ws.Range("D15").Value = "ADDR <8>"

# Auto-fit column D so its width reflects the new, longer text (bestFit).
# (Target authored width is 10.109375 chars; set explicitly since this
# matches the bestFit result Excel computed for the new pin names.)
$ws.Columns.Item(4).ColumnWidth = 9.29

# Move the active selection to D16, matching the new view state
$ws.Range("D16").Select() | Out-Null
